# betexplorer "chile / primera-division / 2023" workbook update
#   - three pairs of already-logged matches had their row order swapped
#     (rows 26/27, 45/46, 63/64) -- only the match-specific columns F:V
#     move, the Indice/pais/torneio/temporada/data_partida columns (A:E)
#     stay put because they already matched up correctly
#   - two freshly scraped matches are appended as rows 211/212

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchColumns {
    param($ws, [int]$row1, [int]$row2)
    # columns F..V = 6..22
    for ($c = 6; $c -le 22; $c++) {
        $cell1 = $ws.Cells.Item($row1, $c)
        $cell2 = $ws.Cells.Item($row2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

Swap-MatchColumns $ws 26 27
Swap-MatchColumns $ws 45 46
Swap-MatchColumns $ws 63 64

function Add-MatchRow {
    param(
        $ws, [int]$row, [int]$formatSourceRow,
        [int]$indice, [double]$dataPartida,
        [string]$home, [int]$homeGols, [string]$away, [int]$awayGols,
        [double]$homeOpenOdds, [string]$homeOpenDt, [double]$homeCloseOdds, [string]$homeCloseDt,
        [double]$drawOpenOdds, [string]$drawOpenDt, [double]$drawCloseOdds, [string]$drawCloseDt,
        [double]$awayOpenOdds, [string]$awayOpenDt, [double]$awayCloseOdds, [string]$awayCloseDt,
        [string]$url
    )

    # Pull formatting (styles) straight from the previous last row so the
    # new rows render identically (bold/border index col, date format, ...)
    $ws.Range("A" + $formatSourceRow + ":V" + $formatSourceRow).Copy()
    $ws.Range("A" + $row + ":V" + $row).PasteSpecial(-4122)

    # "temporada" (col D) is literal text "2023" on every row; copy the
    # value+format from the source row so it keeps its inline-string type
    # instead of being auto-detected as a number.
    $ws.Range("D" + $formatSourceRow).Copy()
    $ws.Range("D" + $row).PasteSpecial(-4104)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($row, 1).Value2 = $indice
    $ws.Cells.Item($row, 2).Value2 = "chile"
    $ws.Cells.Item($row, 3).Value2 = "primera-division"
    $ws.Cells.Item($row, 5).Value2 = $dataPartida
    $ws.Cells.Item($row, 6).Value2 = $home
    $ws.Cells.Item($row, 7).Value2 = $homeGols
    $ws.Cells.Item($row, 8).Value2 = $away
    $ws.Cells.Item($row, 9).Value2 = $awayGols
    $ws.Cells.Item($row, 10).Value2 = $homeOpenOdds
    $ws.Cells.Item($row, 11).Value2 = $homeOpenDt
    $ws.Cells.Item($row, 12).Value2 = $homeCloseOdds
    $ws.Cells.Item($row, 13).Value2 = $homeCloseDt
    $ws.Cells.Item($row, 14).Value2 = $drawOpenOdds
    $ws.Cells.Item($row, 15).Value2 = $drawOpenDt
    $ws.Cells.Item($row, 16).Value2 = $drawCloseOdds
    $ws.Cells.Item($row, 17).Value2 = $drawCloseDt
    $ws.Cells.Item($row, 18).Value2 = $awayOpenOdds
    $ws.Cells.Item($row, 19).Value2 = $awayOpenDt
    $ws.Cells.Item($row, 20).Value2 = $awayCloseOdds
    $ws.Cells.Item($row, 21).Value2 = $awayCloseDt
    $ws.Cells.Item($row, 22).Value2 = $url
}

Add-MatchRow $ws 211 210 210 45242.6875 `
    "Magallanes" 1 "Huachipato" 1 `
    2.8 "07/11/2023 22:12" 3.43 "12/11/2023 16:28" `
    3.32 "07/11/2023 22:12" 3.56 "12/11/2023 16:28" `
    2.49 "07/11/2023 22:12" 2.17 "12/11/2023 16:28" `
    "https://www.betexplorer.com/football/chile/primera-division/magallanes-huachipato/8pELWRah/"

Add-MatchRow $ws 212 211 211 45242.79166666666 `
    "Everton" 1 "Curico Unido" 1 `
    1.47 "06/11/2023 23:12" 1.5 "12/11/2023 18:44" `
    4.53 "06/11/2023 23:12" 4.66 "12/11/2023 18:57" `
    6.84 "06/11/2023 23:12" 6.41 "12/11/2023 18:57" `
    "https://www.betexplorer.com/football/chile/primera-division/everton-curico-unido/ENEHX7pn/"
